$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new text value.
# Values mirror the commit's updated crypto price/volume snapshot.
$updates = @{
    2  = @{ D = "57.426.24"; E = "  -7.03%  " }
    3  = @{ D = "2.883.45";  E = "  -5.38%  " }
    4  = @{ E = "  +0.09%  " }
    5  = @{ D = "550.89";    E = "  -5.94%  " }
    6  = @{ D = "121.23";    E = "  -7.03%  " }
    7  = @{ E = "  +0.17%  " }
    8  = @{ D = "2.875.43";  E = "  -5.63%  " }
    9  = @{ D = "0.489";     E = "  -2.84%  " }
    10 = @{ E = "  -10.84%  " }
    11 = @{ D = "4.77";      E = "  -9.53%  " }
    12 = @{ D = "0.430";     E = "  -2.49%  " }
    13 = @{ D = "0.0000209"; E = "  -11.22%  " }
    14 = @{ D = "31.23";     E = "  -7.35%  " }
    15 = @{ E = "  -0.73%  " }
    16 = @{ D = "3.357.49";  E = "  -5.39%  " }
    17 = @{ D = "2.878.81";  E = "  -5.78%  " }
    18 = @{ D = "57.334.60"; E = "  -7.46%  " }
    19 = @{ D = "6.36";      E = "  -0.60%  " }
    20 = @{ D = "406.94";    E = "  -9.35%  " }
    21 = @{ E = "  -6.40%  " }
    22 = @{ D = "0.649";     E = "  -3.63%  " }
    23 = @{ D = "6.72";      E = "  -8.81%  " }
    24 = @{ D = "12.53";     E = "  -2.70%  " }
    25 = @{ D = "76.50";     E = "  -5.67%  " }
    26 = @{ D = "0.996";     E = "  +0.03%  " }
    27 = @{ E = "  -0.08%  " }
    28 = @{ E = "  -4.65%  " }
    29 = @{ B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "7.10";  E = "  -4.48%  " }
    30 = @{ B = "ImmutableX";  C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";       D = "1.90";  E = "  -6.00%  " }
    31 = @{ E = "  -6.75%  " }
    32 = @{ D = "24.50";     E = "  -5.47%  " }
    33 = @{ D = "0.0947";    E = "  -2.91%  " }
    34 = @{ B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "0.896"; E = "  -8.12%  " }
    35 = @{ B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx";  D = "2.01";  E = "  -13.66%  " }
    36 = @{ E = "  -6.91%  " }
    37 = @{ D = "48.32";     E = "  -4.22%  " }
    38 = @{ D = "8.34";      E = "  +4.85%  " }
    39 = @{ D = "0.0₃0609"; E = "  -12.36%  " }
    40 = @{ D = "0.0342";    E = "  -9.05%  " }
    41 = @{ E = "  -2.82%  " }
    42 = @{ D = "2.596.66";  E = "  -3.61%  " }
    43 = @{ D = "356.25";    E = "  -6.36%  " }
    44 = @{ B = "USDe";      C = "https://coinranking.com/coin/exbfr2U-0+usde-usde";     D = "0.999"; E = "  -0.02%  " }
    45 = @{ B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "2.33";  E = "  -7.21%  " }
    46 = @{ D = "117.43";    E = "  -5.42%  " }
    47 = @{ E = "  -5.67%  " }
    48 = @{ E = "  -3.34%  " }
    49 = @{ E = "  -5.14%  " }
    50 = @{ D = "22.27";     E = "  -7.45%  " }
    51 = @{ D = "1.92";      E = "  -7.57%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $cols[$col]
    }
}
